$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The stray "_GoBack" bookmark currently wraps the inline picture
#    in the result-image paragraph. Word moves this bookmark to track
#    the last editing position, so it needs to be removed from there.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
}

# ------------------------------------------------------------------
# 2) In the "Ket luan" (conclusion) paragraph, the author's cursor
#    ended up right after "tuc la " (i.e. the word "nen " was
#    deleted there). Recreate that: split the run after "tuc la " by
#    dropping a fresh zero-length "_GoBack" bookmark at that point,
#    then strip the leading "nen " word from the remaining text.
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("tức là ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $findRng.End

$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

$rmRng = $d.Range($splitPoint, $d.Content.End)
$rmRng.Find.Execute("nền ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rmRng.Text = ""

Write-Output "Applied ket luan fail update"
